# Generate Report for Handback
# This script swaps the report rows for the two localized files
# (2b34266a-...md and 48e1bd1f-...md) across the Overview, zh-cn and
# de-de worksheets, and refreshes the handback status/timestamps to
# reflect that both files are now "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

# ============================= Overview ==================================
$ws = $wb.Worksheets.Item("Overview")

# Row 2 now describes 2b34266a..., row 3 now describes 48e1bd1f...
$ws.Cells.Item(2,1).Value2 = "2b34266a-9a16-4b77-96b0-a1a636131231.md"
$ws.Cells.Item(2,2).Value2 = "e2e\2b34266a-9a16-4b77-96b0-a1a636131231.md"
$ws.Cells.Item(2,7).Value2 = "2016-08-16 22:47:30"

$ws.Cells.Item(3,1).Value2 = "48e1bd1f-71ff-45e0-ad67-b8d989693882.md"
$ws.Cells.Item(3,2).Value2 = "e2e\48e1bd1f-71ff-45e0-ad67-b8d989693882.md"
$ws.Cells.Item(3,5).Value2 = "Handed back: in sync with en-US"
$ws.Cells.Item(3,6).Value2 = "Handed back: in sync with en-US"
$ws.Cells.Item(3,7).Value2 = "2016-08-16 22:46:45"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\2b34266a-9a16-4b77-96b0-a1a636131231.md"
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\48e1bd1f-71ff-45e0-ad67-b8d989693882.md"
    }
}

# ============================== zh-cn =====================================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Cells.Item(2,1).Value2 = "2b34266a-9a16-4b77-96b0-a1a636131231.md"
$ws.Cells.Item(2,7).Value2 = "2b34266a-9a16-4b77-96b0-a1a636131231.ff5b68a1e153cc32c4141fc6b874641c12a73b71.zh-cn.xlf"
$ws.Cells.Item(2,8).Value2 = "2016-08-16 22:47:24"
$ws.Cells.Item(2,9).Value2 = "2b34266a-9a16-4b77-96b0-a1a636131231.md"
$ws.Cells.Item(2,10).Value2 = "2b34266a-9a16-4b77-96b0-a1a636131231.ff5b68a1e153cc32c4141fc6b874641c12a73b71.zh-cn.xlf"
$ws.Cells.Item(2,11).Value2 = "2016-08-16 22:47:41"

$ws.Cells.Item(3,1).Value2 = "48e1bd1f-71ff-45e0-ad67-b8d989693882.md"
$ws.Cells.Item(3,3).Value2 = "Handed back: in sync with en-US"
$ws.Cells.Item(3,7).Value2 = "48e1bd1f-71ff-45e0-ad67-b8d989693882.fbaa92d89971f24fc89d3e53da47a76236f84fb8.zh-cn.xlf"
$ws.Cells.Item(3,8).Value2 = "2016-08-16 22:46:39"
$ws.Cells.Item(3,9).Value2 = "48e1bd1f-71ff-45e0-ad67-b8d989693882.md"
$ws.Cells.Item(3,10).Value2 = "48e1bd1f-71ff-45e0-ad67-b8d989693882.fbaa92d89971f24fc89d3e53da47a76236f84fb8.zh-cn.xlf"
$ws.Cells.Item(3,16).Value2 = ""

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "2b34266a-9a16-4b77-96b0-a1a636131231.md"
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = "2b34266a-9a16-4b77-96b0-a1a636131231.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "48e1bd1f-71ff-45e0-ad67-b8d989693882.md"
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = "48e1bd1f-71ff-45e0-ad67-b8d989693882.md"
    }
}

$ws.Columns.Item(16).ColumnWidth = 12.86

# ============================== de-de =====================================
$ws = $wb.Worksheets.Item("de-de")

$ws.Cells.Item(2,1).Value2 = "2b34266a-9a16-4b77-96b0-a1a636131231.md"
$ws.Cells.Item(2,7).Value2 = "2b34266a-9a16-4b77-96b0-a1a636131231.ff5b68a1e153cc32c4141fc6b874641c12a73b71.de-de.xlf"
$ws.Cells.Item(2,8).Value2 = "2016-08-16 22:47:30"
$ws.Cells.Item(2,9).Value2 = "2b34266a-9a16-4b77-96b0-a1a636131231.md"
$ws.Cells.Item(2,10).Value2 = "2b34266a-9a16-4b77-96b0-a1a636131231.ff5b68a1e153cc32c4141fc6b874641c12a73b71.de-de.xlf"
$ws.Cells.Item(2,11).Value2 = "2016-08-16 22:47:49"

$ws.Cells.Item(3,1).Value2 = "48e1bd1f-71ff-45e0-ad67-b8d989693882.md"
$ws.Cells.Item(3,3).Value2 = "Handed back: in sync with en-US"
$ws.Cells.Item(3,7).Value2 = "48e1bd1f-71ff-45e0-ad67-b8d989693882.fbaa92d89971f24fc89d3e53da47a76236f84fb8.de-de.xlf"
$ws.Cells.Item(3,8).Value2 = "2016-08-16 22:46:45"
$ws.Cells.Item(3,9).Value2 = "48e1bd1f-71ff-45e0-ad67-b8d989693882.md"
$ws.Cells.Item(3,10).Value2 = "48e1bd1f-71ff-45e0-ad67-b8d989693882.fbaa92d89971f24fc89d3e53da47a76236f84fb8.de-de.xlf"
$ws.Cells.Item(3,16).Value2 = ""

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "2b34266a-9a16-4b77-96b0-a1a636131231.md"
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = "2b34266a-9a16-4b77-96b0-a1a636131231.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "48e1bd1f-71ff-45e0-ad67-b8d989693882.md"
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = "48e1bd1f-71ff-45e0-ad67-b8d989693882.md"
    }
}

$ws.Columns.Item(16).ColumnWidth = 12.86
